$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.121.04'
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").Value = '1.571.53'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.64'
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("E6").Value = '  -3.41%  '
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0608'
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.59'
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.791.66'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.04'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.568.90'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.512'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.11'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '26.119.48'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.26'
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '206.97'
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.76'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.111'
$ws.Range("E28").Value = '  -2.20%  '
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  -1.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.20'
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = '1.278.49'
$ws.Range("E34").Value = '  -0.90%  '
$ws.Range("E35").Value = '  -1.41%  '
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  -1.64%  '
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.09'
$ws.Range("E39").Value = '  -9.17%  '
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("E42").Value = '  -2.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.761'
$ws.Range("E43").Value = '  -2.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.08'
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").Value = '1.705.70'
$ws.Range("E45").Value = '  -2.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.01'
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1000'
$ws.Range("E49").Value = '  -1.80%  '
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("E51").Value = '  -0.35%  '
